# Generate Report for Handoff
#
# A fresh handoff was generated for the files that hadn't yet received their
# own individual "Latest Handoff Datetime" (they still showed the batch
# default value). Update column D ("Latest Handoff Datetime") for those rows
# on both the "zh-cn" and "de-de" status sheets to the new handoff timestamp.

$wb = $excel.ActiveWorkbook

$rows = 4,6,7,8,9,10

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 4).Value = "2016-02-29 04:46:24"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 4).Value = "2016-02-29 04:46:35"
}
